$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells F1:H1, copying the style of the existing header E1
$ws.Range("F1").Value = "KNN_Outliers_MAD"
$ws.Range("G1").Value = "SVM_Outliers_MAD"
$ws.Range("H1").Value = "RF_Outliers_MAD"

$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1:H1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Fill boolean FALSE values for rows 2 through 16 in columns F, G, H
for ($row = 2; $row -le 16; $row++) {
    $ws.Cells.Item($row, 6).Value = $false
    $ws.Cells.Item($row, 7).Value = $false
    $ws.Cells.Item($row, 8).Value = $false
}
